# "Updated CVDs for the month"
# Applies the monthly CVD (ytd column "E") refresh, plus the knock-on
# recompute of the quarterly/annual run-rate columns (O:W) that are derived
# from the CVD, across the affected site tabs.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Cassville Missouri
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Cassville Missouri")
$ws.Range("E2").Value2 = 0.0621
$ws.Range("E3").Value2 = 0.0621
$ws.Range("E4").Value2 = 0.0621
$ws.Range("O4:W4").Value2 = 0

# ---------------------------------------------------------------------
# Tipp City Ohio
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Tipp City Ohio")
$ws.Range("E2").Value2 = 0.2941
$ws.Range("E3").Value2 = 0.2941
$ws.Range("E4").Value2 = 0.2941
$ws.Range("O4:W4").Value2 = 0

# ---------------------------------------------------------------------
# Milwaukee Pmc Hq Wisconsin
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Milwaukee Pmc Hq Wisconsin")
$ws.Range("O3").ClearContents()

# ---------------------------------------------------------------------
# Fort Wayne Indiana
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Fort Wayne Indiana")
$ws.Range("E2").Value2 = 0.0836
$ws.Range("E3").Value2 = 0.0836
$ws.Range("E4").Value2 = 0.0836
$ws.Range("O4:W4").Value2 = 0

# ---------------------------------------------------------------------
# Hyderabad India
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Hyderabad India")
$ws.Range("E2").Value2 = 0.0509
$ws.Range("E3").Value2 = 0.0509
$ws.Range("E4").Value2 = 0.0509
$ws.Range("O4:W4").Value2 = 0

# ---------------------------------------------------------------------
# Lincoln Missouri
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Lincoln Missouri")
$ws.Range("O4:W4").Value2 = 0
$ws.Range("O7").ClearContents()

$ws.Range("E8").Value2 = 0.1112
$ws.Range("E9").Value2 = 0.1112
$ws.Range("E10").Value2 = 0.1112

$ws.Range("K10").Value2 = 0.0444
$ws.Range("L10").Value2 = 0.0345
$ws.Range("M10").Value2 = 0.0116
$ws.Range("N10").Value2 = 0.0912
$ws.Range("O10").Value2 = 0
$ws.Range("P10").Value2 = 0.0158833333333333
$ws.Range("Q10").Value2 = 0.0158833333333333
$ws.Range("R10").Value2 = 0.04765
$ws.Range("S10").Value2 = 0.0158833333333333
$ws.Range("T10").Value2 = 0.0158833333333333
$ws.Range("U10").Value2 = 0.0158833333333333
$ws.Range("V10").Value2 = 0.04765
$ws.Range("W10").Value2 = 0.1906

# ---------------------------------------------------------------------
# Piedras Negras Jakel Mexico
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Piedras Negras Jakel Mexico")
$ws.Range("E2").Value2 = 0.0526
$ws.Range("E3").Value2 = 0.0526
$ws.Range("E4").Value2 = 0.0526
$ws.Range("O4:W4").Value2 = 0
$ws.Range("O5").ClearContents()

$ws.Range("E6").Value2 = 0.1728
$ws.Range("E7").Value2 = 0.1728
$ws.Range("E8").Value2 = 0.1728

$ws.Range("G8").Value2 = 0.0315
$ws.Range("I8").Value2 = 0.0111
$ws.Range("J8").Value2 = 0.0554
$ws.Range("K8").Value2 = 0.033
$ws.Range("L8").Value2 = 0.0236
$ws.Range("M8").Value2 = 0.0246
$ws.Range("N8").Value2 = 0.0811
$ws.Range("O8").Value2 = 0.0373
$ws.Range("P8").Value2 = 0.0246833333333333
$ws.Range("Q8").Value2 = 0.0246833333333333
$ws.Range("R8").Value2 = 0.07405
$ws.Range("S8").Value2 = 0.0246833333333333
$ws.Range("T8").Value2 = 0.0246833333333333
$ws.Range("U8").Value2 = 0.0246833333333333
$ws.Range("V8").Value2 = 0.07405
$ws.Range("W8").Value2 = 0.2962
